# Rough Draft BDM Instructions
# Update the three BDM question prompts (Simon / Kriti / Borys & Sabi) to
# reference "the lottery (i.e case 2)" instead of "the probability of leaving
# the experiment right now", tweak wording, and update the Borys & Sabi
# question's final clause. Also move the active cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "I used the BDM procedure to determine how much Simon values an ice-cream in terms of the lottery (i.e case 2). Simon (truthfully) stated that his switch-point is 60%. From this, we can infer:"

$ws.Range("A6").Value = "I used the BDM procedure to determine how much Kriti values a mug in terms of the lottery (i.e case 2). Suppose Kriti states that her switch-point for 1 mug is 100%,and her switch point for 2 mugs is also 100%. Kriti also (truthfully) states that she prefers 2 mugs to 1 mug. What can we infer?"

$ws.Range("A9").Value = "I used the BDM procedure to determine how much Borys and Sabi value a coffee in terms of the probability of leaving the experiment right now. Borys stated his switch-point was 25%, while Sabi stated his switch point was 50%. Both Borys and Sabi prefer leaving the experiment right now to reading the boring information. What can we infer?"

# The shorter question text in row 5 now only needs a 75pt row instead of 90pt.
$ws.Rows.Item(5).RowHeight = 75

# Move the selection from G4 to A6, matching the author's final cursor position.
$ws.Range("A6").Select()
